$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Narrow column G (Status/Time) width from 18 to 8 characters.
# Excel's ColumnWidth COM property is expressed in the default font's
# character units and gets re-derived into the OOXML "width" value, so we
# use the value that yields a stored width of exactly 8.
$ws.Columns.Item(7).ColumnWidth = 7.16666665

# Update game status from "17:39 - 2nd Half" to "Final" for the finished games
$rows = @(7, 12, 13, 16, 17, 21, 25, 27, 28, 32, 33, 35, 37, 66, 67, 73, 76, 79)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "Final"
}
